$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 1.065782333333333
$ws.Range("N2").Value = 3.197347
$ws.Range("O2").Value = 0.1896668697192897
$ws.Range("P2").Value = 0.1896668697192897
$ws.Range("Q2").Value = 214.6909189653463
$ws.Range("R2").Value = 1932.218270688116
$ws.Range("S2").Value = 0.09167548589969532
$ws.Range("T2").Value = 0.09167548589969535
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.009797007290259483
$ws.Range("P3").Value = 0.009797007290259485
$ws.Range("Q3").Value = 11.08959356670445
$ws.Range("R3").Value = 99.80634210034
$ws.Range("S3").Value = 0.00473538370210183
$ws.Range("T3").Value = 0.004735383702101831
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 1.641159333333333
$ws.Range("N4").Value = 4.923477999999999
$ws.Range("O4").Value = 0.2920610932725753
$ws.Range("P4").Value = 0.2920610932725753
$ws.Range("Q4").Value = 330.5947137816649
$ws.Range("R4").Value = 2975.352424034984
$ws.Range("S4").Value = 0.1411677362408459
$ws.Range("T4").Value = 0.141167736240846
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 2.857239666666667
$ws.Range("N5").Value = 8.571719
$ws.Range("O5").Value = 0.5084750297178755
$ws.Range("P5").Value = 0.5084750297178756
$ws.Range("Q5").Value = 575.5616231903258
$ws.Range("R5").Value = 5180.054608712932
$ws.Range("S5").Value = 0.2457714174659962
$ws.Range("T5").Value = 0.2457714174659962
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 1.065782333333333
$ws.Range("N6").Value = 3.197347
$ws.Range("O6").Value = 0.1896668697192897
$ws.Range("P6").Value = 0.1896668697192897
$ws.Range("Q6").Value = 69.72067155495756
$ws.Range("R6").Value = 627.4860439946181
$ws.Range("S6").Value = 0.02977152677373138
$ws.Range("T6").Value = 0.02977152677373139
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("O7").Value = 0.009797007290259483
$ws.Range("P7").Value = 0.009797007290259485
$ws.Range("S7").Value = 0.001537811349320423
$ws.Range("T7").Value = 0.001537811349320423
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("M8").Value = 1.641159333333333
$ws.Range("N8").Value = 4.923477999999999
$ws.Range("O8").Value = 0.2920610932725753
$ws.Range("P8").Value = 0.2920610932725753
$ws.Range("Q8").Value = 107.3603185847702
$ws.Range("R8").Value = 966.2428672629319
$ws.Range("S8").Value = 0.045844087956946
$ws.Range("T8").Value = 0.04584408795694601
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 2.857239666666667
$ws.Range("N9").Value = 8.571719
$ws.Range("O9").Value = 0.5084750297178755
$ws.Range("P9").Value = 0.5084750297178756
$ws.Range("Q9").Value = 186.9130892143985
$ws.Range("R9").Value = 1682.217802929586
$ws.Range("S9").Value = 0.07981403385538134
$ws.Range("T9").Value = 0.07981403385538137
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 1.065782333333333
$ws.Range("N10").Value = 3.197347
$ws.Range("O10").Value = 0.1896668697192897
$ws.Range("P10").Value = 0.1896668697192897
$ws.Range("Q10").Value = 64.41038621086979
$ws.Range("R10").Value = 579.6934758978281
$ws.Range("S10").Value = 0.02750397399818125
$ws.Range("T10").Value = 0.02750397399818126
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("O11").Value = 0.009797007290259483
$ws.Range("P11").Value = 0.009797007290259485
$ws.Range("Q11").Value = 3.327038740135556
$ws.Range("R11").Value = 29.94334866122
$ws.Range("S11").Value = 0.001420683718617224
$ws.Range("T11").Value = 0.001420683718617224
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("M12").Value = 1.641159333333333
$ws.Range("N12").Value = 4.923477999999999
$ws.Range("O12").Value = 0.2920610932725753
$ws.Range("P12").Value = 0.2920610932725753
$ws.Range("Q12").Value = 99.18320391271911
$ws.Range("R12").Value = 892.648835214472
$ws.Range("S12").Value = 0.04235236616251456
$ws.Range("T12").Value = 0.04235236616251457
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 2.857239666666667
$ws.Range("N13").Value = 8.571719
$ws.Range("O13").Value = 0.5084750297178755
$ws.Range("P13").Value = 0.5084750297178756
$ws.Range("Q13").Value = 172.6768259063062
$ws.Range("R13").Value = 1554.091433156756
$ws.Range("S13").Value = 0.07373498606679732
$ws.Range("T13").Value = 0.07373498606679735
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 1.065782333333333
$ws.Range("N14").Value = 3.197347
$ws.Range("O14").Value = 0.1896668697192897
$ws.Range("P14").Value = 0.1896668697192897
$ws.Range("Q14").Value = 95.35079375043057
$ws.Range("R14").Value = 858.1571437538752
$ws.Range("S14").Value = 0.04071588304768173
$ws.Range("T14").Value = 0.04071588304768175
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("O15").Value = 0.009797007290259483
$ws.Range("P15").Value = 0.009797007290259485
$ws.Range("Q15").Value = 4.925227177986112
$ws.Range("R15").Value = 44.327044601875
$ws.Range("S15").Value = 0.002103128520220006
$ws.Range("T15").Value = 0.002103128520220007
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("M16").Value = 1.641159333333333
$ws.Range("N16").Value = 4.923477999999999
$ws.Range("O16").Value = 0.2920610932725753
$ws.Range("P16").Value = 0.2920610932725753
$ws.Range("Q16").Value = 146.8272087179722
$ws.Range("R16").Value = 1321.44487846175
$ws.Range("S16").Value = 0.06269690291226881
$ws.Range("T16").Value = 0.06269690291226883
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 2.857239666666667
$ws.Range("N17").Value = 8.571719
$ws.Range("O17").Value = 0.5084750297178755
$ws.Range("P17").Value = 0.5084750297178756
$ws.Range("Q17").Value = 255.6244944498195
$ws.Range("R17").Value = 2300.620450048375
$ws.Range("S17").Value = 0.1091545923297007
$ws.Range("T17").Value = 0.1091545923297007
